$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.45 = 5464.65 pesos`n✅ 5464.65 pesos = 1.44 = 783.13 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 690
$ws2.Range("O10").Value = 3770.61
$ws2.Range("N12").Value = 3803
$ws2.Range("O12").Value = 545
